$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 (year 2018): add the actual expenses formula, matching the
# pattern used in the rows above it (hard-coded known totals).
$ws.Range("E8").Formula = "=-2234.42-C8-D8"

# G8 becomes a plain (known/actual) rent value rather than the projected
# growth formula.
$ws.Range("G8").Value = 2100

# H8: newly entered known expense value.
$ws.Range("H8").Value = -1330.14

# Leave the selection where the user ended up after editing E8 (Enter
# moves the active cell down one row).
$ws.Range("E9").Select()
